$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.702000000000001
$ws.Range("D4").Value = -7.945
$ws.Range("C7").Value = -12.983
$ws.Range("A8").Value = -22.231
$ws.Range("A10").Value = -22.148
$ws.Range("E10").Value = 16.451
$ws.Range("D11").Value = -7.384
$ws.Range("A12").Value = -21.663
$ws.Range("E12").Value = 17.578
$ws.Range("E13").Value = 16.416
$ws.Range("C14").Value = -13.291
$ws.Range("D14").Value = -8.196
$ws.Range("E14").Value = 16.583
$ws.Range("C15").Value = -13.728
$ws.Range("A18").Value = -21.862
$ws.Range("C18").Value = -12.226
$ws.Range("D18").Value = -7.768000000000001
$ws.Range("D19").Value = -8.254999999999999
$ws.Range("C20").Value = -12.684
$ws.Range("D21").Value = -8.416
$ws.Range("A25").Value = -21.953
$ws.Range("D27").Value = -8.068000000000001
$ws.Range("C29").Value = -12.591
$ws.Range("E29").Value = 16.747
$ws.Range("C30").Value = -12.392
$ws.Range("C31").Value = -12.832
$ws.Range("D31").Value = -7.858
$ws.Range("E32").Value = 17.057
$ws.Range("C35").Value = -12.532
$ws.Range("E35").Value = 16.525
$ws.Range("A37").Value = -20.044
$ws.Range("D38").Value = -7.759
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.304
$ws.Range("E43").Value = 16.789
$ws.Range("C44").Value = -12.993
$ws.Range("D44").Value = -7.592000000000001
$ws.Range("D47").Value = -7.395
$ws.Range("E48").Value = 17.176
$ws.Range("E49").Value = 16.449
$ws.Range("C50").Value = -12.784
$ws.Range("E50").Value = 16.365
$ws.Range("E51").Value = 16.757
$ws.Range("C54").Value = -12.524
$ws.Range("A55").Value = -21.946
$ws.Range("D56").Value = -8.19
$ws.Range("E56").Value = 16.237
$ws.Range("D58").Value = -8.401
$ws.Range("E61").Value = 16.461
$ws.Range("D65").Value = -7.928999999999999
$ws.Range("A68").Value = -21.602
$ws.Range("C68").Value = -11.281
$ws.Range("E69").Value = 17.32
$ws.Range("E71").Value = 17.31
$ws.Range("D73").Value = -8.193000000000001
$ws.Range("C76").Value = -13.314
$ws.Range("A77").Value = -20.262
$ws.Range("A78").Value = -20.099
$ws.Range("A79").Value = -21.449
$ws.Range("E79").Value = 16.748
$ws.Range("A80").Value = -20.228
$ws.Range("A81").Value = -21.698
$ws.Range("E81").Value = 16.67
$ws.Range("A82").Value = -22.124
$ws.Range("A84").Value = -22.177
$ws.Range("C87").Value = -13.215
$ws.Range("C88").Value = -13.013
$ws.Range("D90").Value = -7.43
$ws.Range("C92").Value = -11.733
$ws.Range("D92").Value = -7.551
$ws.Range("E92").Value = 16.724
$ws.Range("D94").Value = -7.355
$ws.Range("D95").Value = -7.836
$ws.Range("C96").Value = -12.665
$ws.Range("C98").Value = -13.649
$ws.Range("A101").Value = -20.761
$ws.Range("C101").Value = -12.592
$ws.Range("D101").Value = -7.830000000000001
$ws.Range("A102").Value = -19.812
$ws.Range("C102").Value = -12.567
